$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (also updates the "Through 2022-03-03" -> "Through 2022-03-04" in workbook.xml)
$ws.Name = "Through 2022-03-04"

# Update the header label in I1 (shared string "2022 (through 03-03)" -> "2022 (through 03-04)")
$ws.Range("I1").Value = "2022 (through 03-04)"

# Update the March (row 3) and April (row 4) values in column I
$ws.Range("I3").Value = 142
$ws.Range("I4").Value = 20

# Update the Total (row 14) value in column I
$ws.Range("I14").Value = 321
